# Apply updated values to column F ("dSF") on Sheet1, per repull/recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 5
    3  = -6
    7  = -2
    9  = -3
    10 = -4
    12 = -5
    13 = -1
    14 = -4
    17 = -3
    20 = 6
    21 = 6
    24 = 2
    25 = -3
    26 = 5
    35 = -5
    37 = 0
    39 = 1
    42 = -4
    47 = -3
    49 = 1
    53 = -2
    58 = -3
    60 = -2
    64 = 2
    66 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
